$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 updates
$ws.Range("G11").Value = 1.67
$ws.Range("I11").Value = 5.5
$ws.Range("J11").Value = 2.3
$ws.Range("U11").Value = 2.1
$ws.Range("V11").Value = 1.67
$ws.Range("X11").Value = 7
$ws.Range("Z11").Value = 12
$ws.Range("AD11").Value = 7
$ws.Range("AJ11").Value = 19
$ws.Range("AN11").Value = 3.5

# Row 12 updates
$ws.Range("G12").Value = 4.65
$ws.Range("H12").Value = 3.7
$ws.Range("I12").Value = 1.65
$ws.Range("J12").Value = 4.75
$ws.Range("K12").Value = 2.22
$ws.Range("L12").Value = 2.18
$ws.Range("M12").Value = 1.03
$ws.Range("N12").Value = 11
$ws.Range("O12").Value = 1.22
$ws.Range("S12").Value = 1.36
$ws.Range("T12").Value = 2.94
$ws.Range("U12").Value = 1.65
$ws.Range("V12").Value = 1.98
$ws.Range("W12").Value = 14
$ws.Range("X12").Value = 28
$ws.Range("Y12").Value = 15
$ws.Range("Z12").Value = 80
$ws.Range("AA12").Value = 45
$ws.Range("AB12").Value = 45
$ws.Range("AC12").Value = 11.75
$ws.Range("AD12").Value = 7.3
$ws.Range("AE12").Value = 14.5
$ws.Range("AF12").Value = 60
$ws.Range("AG12").Value = 400
$ws.Range("AH12").Value = 7.7
$ws.Range("AI12").Value = 8.5
$ws.Range("AJ12").Value = 8
$ws.Range("AK12").Value = 13
$ws.Range("AL12").Value = 12.5
$ws.Range("AM12").Value = 23
$ws.Range("AN12").Value = 6.4
$ws.Range("AO12").Value = 25
$ws.Range("AP12").Value = 29
$ws.Range("AQ12").Value = 150
$ws.Range("AR12").Value = 175
$ws.Range("AS12").Value = 350
$ws.Range("AT12").Value = 2.9
$ws.Range("AU12").Value = 7.1
$ws.Range("AV12").Value = 60
$ws.Range("AW12").Value = 3.55
$ws.Range("AX12").Value = 7.9
$ws.Range("AY12").Value = 16
$ws.Range("AZ12").Value = 25
$ws.Range("BA12").Value = 50
